$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing Caso "5507" / "CONGRESO 1927" (originally row 19).
$ws.Rows(19).Delete()

# After the row above shifts everything up by one, the row containing Caso
# "6321" / "ARCOS 2739" (originally row 48) is now at row 47. Remove it too.
$ws.Rows(47).Delete()
